$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 21: refreshed "dct:modified" timestamp (new .ttl generation run) ---
$ws.Range("B21").Value = "2023-09-13T15:57:50+00:00"

# --- Row 23: header mapping row C:Z re-ordered / updated ---
# skos:altLabel moves from G to C, skos:notation is newly inserted at D,
# skos:broader moves from C to G, rdf:type moves from D to V, and the
# trailing dct: columns each shift one place to the right (W..Z).
$ws.Range("C23").Value = 'skos:altLabel(separator=",")'
$ws.Range("D23").Value = "skos:notation"
$ws.Range("G23").Value = 'skos:broader(separator=",")'
$ws.Range("V23").Value = "rdf:type"
$ws.Range("W23").Value = "dct:modified^^xsd:date"
$ws.Range("X23").Value = "dct:created^^xsd:date"
$ws.Range("Y23").Value = 'dct:creator(separator=",")'
$ws.Range("Z23").Value = 'dct:contributor(separator=",")'

# --- Row 25 ("variable" term): add the missing notation ---
$ws.Range("D25").Value = "var"

# --- Row 27 ("emerging" term): C held the broader-concept link, now holds
#     the term's own notation, and the broader-concept link moves to G ---
$ws.Range("C27").Value = "new"
$ws.Range("G27").Value = "vocab:1002"

# --- Row 28 ("developing" term): same shuffle as row 27 ---
$ws.Range("C28").Value = "intermediate"
$ws.Range("G28").Value = "vocab:1002"

# --- Row 29 ("mature" term): same shuffle, but no notation replaces C ---
# (cleared to blank text, mirroring the sheet's other empty text cells,
# rather than removed outright)
$ws.Range("C29").Value = "'"
$ws.Range("C29").Style = "Normal"
$ws.Range("G29").Value = "vocab:1002"

# --- Row 30 ("hasMaturityLevel" term): rdf:type value moves from D to V ---
$ws.Range("D30").Value = "'"
$ws.Range("D30").Style = "Normal"
$ws.Range("V30").Value = "owl:ObjectProperty"

# --- New column AO: the freshly regenerated sheet has one extra (blank)
#     trailing text column across every existing row, matching the blank
#     inline-string cells already present in columns D..AN. Forcing a
#     text value (via a leading apostrophe) then emptying the resulting
#     quote-prefix style keeps the cells text-typed-but-blank, exactly
#     like their neighbours, instead of leaving them as unused numeric
#     cells. ---
for ($r = 1; $r -le 30; $r++) {
    $cell = $ws.Cells.Item($r, 41)
    $cell.Value = "'"
}
for ($r = 1; $r -le 30; $r++) {
    $ws.Cells.Item($r, 41).Style = "Normal"
}

Write-Host "edit applied"
